# Update COVID-19 "paises" workbook to the later snapshot (20:22 update).
# Commit message: "Update countries & provincias Spain"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp update -------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 20:22"

# Each tuple: Row, Country, Total, NewCases, Active, Recovered, Critical, DeathsToday, Deaths
$rows = @(
    # Simple numeric refreshes (no row reordering)
    ,(4,   "Estados Unidos",       577332, 17032, 33907, 520348, 12565, 972, 23077)
    ,(8,   "Alemania",             128166, 312,   64300, 60825,  4895,  19,  3041)
    ,(15,  "Suiza",                25688,  273,   13700, 10850,  386,   32,  1138)
    ,(16,  "Canada",               25548,  1165,  7659,  17122,  557,   50,  767)
    ,(20,  "Austria",              14040,  95,    7343,  6329,   239,   18,  368)
    ,(74,  "Bosnia y Herzegovina", 1037,   28,    206,   792,    4,     0,   39)
    ,(76,  "Uzbekistan",           998,    133,   85,    909,    8,     0,   4)
    ,(113, "Georgia",              272,    15,    68,    201,    6,     0,   3)

    # Peru overtakes Ecuador & Chile (rows 26-28)
    ,(26, "Peru",    9784, 2265, 1798, 7793, 134, 0,  193)
    ,(27, "Ecuador", 7529, 63,   597,  6577, 121, 22, 355)
    ,(28, "Chile",   7525, 312,  2367, 5076, 387, 2,  82)

    # Sudafrica overtakes Argentina (rows 54-55)
    ,(54, "Sudafrica", 2272, 99, 410, 1835, 7,  2, 27)
    ,(55, "Argentina", 2208, 66, 515, 1596, 83, 7, 97)

    # Tunez overtakes Cuba & Crucero (rows 82-84)
    ,(82, "Tunez",   726, 19, 43,  649, 89, 3, 34)
    ,(83, "Cuba",    726, 57, 121, 584, 11, 3, 21)
    ,(84, "Crucero", 712, 0,  619, 82,  10, 0, 11)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
